$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Employees")

# The previous "add employee" attempt left row 6 with stale data
# (EMP ID E0112 / p12) instead of creating a fresh row; correct those
# three text cells to the real new hire (E0114 / p14).
$ws.Range("B6").Value = "E0114"
$ws.Range("C6").Value = "p14"
$ws.Range("G6").Value = "p14@gmail.com"

# Now actually add the next employee as a new row 7, matching the layout
# and formatting already used by the other data rows.
$ws.Range("A6:K6").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "E0115"
$ws.Range("C7").Value = "p15"
$ws.Range("D7").Value2 = 44922
$ws.Range("E7").Value2 = 35954
$ws.Range("F7").Value = "Java Developer"
$ws.Range("G7").Value = "p15@gmail.com"
$ws.Range("H7").Value = "Female"
$ws.Range("I7").Value2 = 1234543268
$ws.Range("J7").Value = "Pune"
$ws.Range("K7").Value2 = 2

$ws.Range("A1").Select()
$ws.Range("F10").Select()
